$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a cell value split across runs with per-run font overrides.
# $segments is an array of 2-element arrays: @(text, fontNameOrNull)
# A $null font leaves that run with no explicit font override (inherits the
# cell/default formatting), a non-null font sets Name+Size(10) on that run.
# ---------------------------------------------------------------------------
function Set-RichText {
    param($cell, $segments)

    $full = ""
    foreach ($seg in $segments) {
        $full += $seg[0]
    }
    $cell.Value = $full

    $pos = 1
    foreach ($seg in $segments) {
        $text = $seg[0]
        $font = $seg[1]
        $len = $text.Length
        if ($len -gt 0 -and $font -ne $null) {
            $r = $cell.Characters($pos, $len)
            $r.Font.Name = $font
            $r.Font.Size = 10
        }
        $pos += $len
    }
}

# ---------------------------------------------------------------------------
# New row 25 / row 26 : inherit formatting from row 24 (the last existing
# entry), same pattern used throughout the sheet (date / hashtag / time /
# content / reference columns).
# ---------------------------------------------------------------------------
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A24:D24").Copy()
$ws.Range("A26:D26").PasteSpecial(-4122)

# Set B26 ("Ajax") first so it claims the shared-string slot right after the
# pre-existing strings, matching the authoring order of the real edit.
$ws.Cells.Item(26, 2).Value = "Ajax"

# D25 : "나는 오늘 버블링과 캡처링, preventDefalut(), Uploadprofile하는것을 배움"
$d25Segments = @(
    ,@("나는", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("오늘", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("버블링과", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("캡처링", "맑은 고딕")
    ,@(", preventDefalut(), Uploadprofile", "Arial")
    ,@("하는것을", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("배움", "맑은 고딕")
)
Set-RichText $ws.Cells.Item(25, 4) $d25Segments

# E25 : full HTML/JS/CSS reference code block (two bubbling/capturing demos).
$e25Text = @'
<!DOCTYPE html>
<html lang="en">
<head>
  <script src="https://code.jquery.com/jquery-3.7.1.min.js"></script>
  <link rel="stylesheet" href="style.css">
</head>
<body>
  <div id ="di">
<a class="link" href="#none">link</a>
</div>
  <script src="cl.js"></script>
</body>
</html>
let n = document.getElementById('di'); // 상위요소
let m = document.querySelector('#di a'); // 하위요소
//상위 요소로 전달되는 버블링,
//하위요소로 전달되는 캡처링
n.onclick=function(){
    this.style.borderBlockColor='red';
}
m.onclick=function(e){
    e.stopPropagation(); //버블링, 캡처링 막는 함수
    this.style.borderBlockColor='blue';
}
body{
    margin: 20px;
}
#di{
    border: 1px solid black;
    padding: 30px;
}
#di a{
    text-decoration: none;
    border: 1px solid black;
    padding: 10px;
}
<!DOCTYPE html>
<html lang="en">
<head>
  <script src="https://code.jquery.com/jquery-3.7.1.min.js"></script>
  <link rel="stylesheet" href="style.css">
</head>
<body>
  <div id ="di">
<a class="link" href="http://naver.com">네입버</a>
</div>
  <script src="cl.js"></script>
</body>
</html>
let n = document.querySelector('.link');
n.addEventListener('click',function(e){
    e.preventDefault();
let m = this.getAttribute('href');
console.log(m);
let z = confirm('네이버로 이동할려?');
if(z==true){
    location.href=m;
}
});
'@
$ws.Cells.Item(25, 5).Value = $e25Text

# D26 : "Tomcat 설치 및 서버구성하는 법, 간단한 동작 실행해봄"
$d26Segments = @(
    ,@("Tomcat ", $null)
    ,@("설치", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("및", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("서버구성하는", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("법", "맑은 고딕")
    ,@(", ", "Arial")
    ,@("간단한", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("동작", "맑은 고딕")
    ,@(" ", "Arial")
    ,@("실행해봄", "맑은 고딕")
)
Set-RichText $ws.Cells.Item(26, 4) $d26Segments

# Remaining plain columns (reuse existing shared strings: Jquery/오전/오후 etc.)
$ws.Cells.Item(25, 1).Value = 45394
$ws.Cells.Item(25, 2).Value = "Jquery"
$ws.Cells.Item(25, 3).Value = "오전"

$ws.Cells.Item(26, 1).Value = 45394
$ws.Cells.Item(26, 3).Value = "오후"

# Pin both new rows back to the sheet's standard row height (customHeight) -
# content changes above can otherwise trigger auto-fit on the long E25 text.
$ws.Rows(25).RowHeight = $ws.Rows(24).RowHeight
$ws.Rows(26).RowHeight = $ws.Rows(24).RowHeight

# ---------------------------------------------------------------------------
# View state: selection moves to E27, and the sheet scrolls back to the top.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("E27").Select()
